$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4 (Union Homes REIT) is removed entirely from the dataset ---
$ws.Rows.Item(4).Delete()

# --- Row 2 (Union Homes REIT, now re-indexed) ---
# B2 label changes from "2" to "1" but must remain a text value (quote-prefixed
# so Excel doesn't silently coerce it to a number), matching the original
# inlineStr cell type.
$ws.Range("B2").Value = "'1"
$ws.Range("B2").Style = "Normal"

$ws.Range("G2").Value  = 1.010830324909747
$ws.Range("H2").Value  = 1.010830324909747
$ws.Range("I2").Value  = 0.7021660649819494
$ws.Range("J2").Value  = 0.7021660649819494
$ws.Range("K2").Value  = 4.1
$ws.Range("L2").Value  = 0.740072202166065
$ws.Range("M2").Value  = 0
$ws.Range("N2").Value  = 0
$ws.Range("O2").Value  = 0
$ws.Range("P2").Value  = 0
$ws.Range("Q2").Value  = 0
$ws.Range("R2").Value  = 0
$ws.Range("S2").Value  = 0
# T2 is dropped entirely (no longer present in the data)
$ws.Range("T2").ClearContents()
$ws.Range("U2").Value  = 4.07
$ws.Range("V2").Value  = 0.1057142857142857
$ws.Range("W2").Value  = 0.04357066950053135
$ws.Range("X2").Value  = 0.07539266199954733
$ws.Range("Y2").Value  = -0.03182199249901598
$ws.Range("Z2").Value  = 0.05951869359690589
$ws.Range("AA2").Value = 0.04179200687580576
$ws.Range("AB2").Value = 0.07539266199954733
$ws.Range("AC2").Value = -0.03360065512374157
$ws.Range("AG2").Value = -4.07
$ws.Range("AJ2").Value = -0.1182108626198083
$ws.Range("AK2").Value = -0.05169566874126763
$ws.Range("AP2").Value = -0.9667458432304039

# --- Row 3 (UPDC REIT) ---
$ws.Range("G3").Value  = 1.010830324909747
$ws.Range("H3").Value  = 1.010830324909747
$ws.Range("I3").Value  = 0.7021660649819494
$ws.Range("J3").Value  = 0.7021660649819494
$ws.Range("K3").Value  = 4.1
$ws.Range("L3").Value  = 0.740072202166065
$ws.Range("M3").Value  = 0
$ws.Range("N3").Value  = 0
$ws.Range("O3").Value  = 0
$ws.Range("P3").Value  = 0
$ws.Range("Q3").Value  = 0
$ws.Range("R3").Value  = 0
$ws.Range("S3").Value  = 0
# T3 is dropped entirely (no longer present in the data)
$ws.Range("T3").ClearContents()
$ws.Range("U3").Value  = 4.07
$ws.Range("V3").Value  = 0.1057142857142857
$ws.Range("W3").Value  = 0.04357066950053135
$ws.Range("X3").Value  = 0.07539266199954733
$ws.Range("Y3").Value  = -0.03182199249901598
$ws.Range("Z3").Value  = 0.05951869359690589
$ws.Range("AA3").Value = 0.04179200687580576
$ws.Range("AB3").Value = 0.07539266199954733
$ws.Range("AC3").Value = -0.03360065512374157
$ws.Range("AG3").Value = -4.07
$ws.Range("AJ3").Value = -0.1182108626198083
$ws.Range("AK3").Value = -0.05169566874126763
$ws.Range("AP3").Value = -0.9667458432304039
